# Update "想去人数" (interested-attendee count) figures in the F column
# across the three data sheets, as published to gh-pages at 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 5576
$ws1.Range("F4").Value  = 7662
$ws1.Range("F11").Value = 4432
$ws1.Range("F12").Value = 1795
$ws1.Range("F15").Value = 2994
$ws1.Range("F20").Value = 474
$ws1.Range("F23").Value = 119
$ws1.Range("F24").Value = 1725
$ws1.Range("F26").Value = 104
$ws1.Range("F27").Value = 1452
$ws1.Range("F33").Value = 21
$ws1.Range("F35").Value = 108
$ws1.Range("F37").Value = 3126
$ws1.Range("F38").Value = 718
$ws1.Range("F40").Value = 159
$ws1.Range("F42").Value = 888

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 23

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 5576
$ws4.Range("F4").Value  = 7662
$ws4.Range("F11").Value = 4432
$ws4.Range("F12").Value = 1795
$ws4.Range("F15").Value = 2994
$ws4.Range("F20").Value = 474
$ws4.Range("F24").Value = 119
$ws4.Range("F25").Value = 1725
$ws4.Range("F27").Value = 104
$ws4.Range("F28").Value = 1452
$ws4.Range("F34").Value = 21
$ws4.Range("F36").Value = 108
$ws4.Range("F38").Value = 3126
$ws4.Range("F39").Value = 23
$ws4.Range("F40").Value = 718
$ws4.Range("F42").Value = 159
$ws4.Range("F44").Value = 888
